$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values B2 and B3 from NAN text to real numbers
$ws.Range("B2").Value = 100
$ws.Range("B3").Value = 73.361080202582

# Update formulas in B4 and B5 to divide by (12*100) instead of (12*500)
$ws.Range("B4").Formula = "=B2/(12*100)"
$ws.Range("B5").Formula = "=B3/(12*100)"

# Update B6 static value
$ws.Range("B6").Value = 0.081194294852525
